$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.998.10"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "1.635.52"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.993"
$ws.Range("E4").Value = "  -0.82%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.95"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.74"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.25"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").Value = "1.858.51"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").Value = "1.618.62"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.18"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "25.972.54"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.993"
$ws.Range("E19").Value = "  -0.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.43"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.78"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.35"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E24").Value = "  -1.00%  "
$ws.Range("E25").Value = "  -2.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.96"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.89"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0498"
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.904"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "1.136.91"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.993"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.58"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.10"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.789"
$ws.Range("E44").Value = "  -2.17%  "
$ws.Range("D45").Value = "1.769.39"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("E46").Value = "  -2.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.76"
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("E48").Value = "  +2.25%  "
$ws.Range("E49").Value = "  +4.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.69"
$ws.Range("E50").Value = "  +1.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.416"
$ws.Range("E51").Value = "  -0.13%  "
